# office.xlsx test fixture update:
# - drop the "id" column (A) from the office listing
# - re-type the header row with capitalised labels
# - keep the numeric-looking code columns as quote-prefixed text
# - widen column B (office name) and move the selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 2): "id" column removed, remaining headers capitalised
$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "OfficeTypeCode"
$ws.Range("D2").Value = "ParentOfficeCode"

# Row 3 ("Красный районный участок" office): shift code/name/type/parent one column left
$ws.Range("A3").Value = "'231100"
$ws.Range("B3").Value = "Красный районный участок"
$ws.Range("C3").Value = "РайУч"
$ws.Range("D3").Value = "'231000"

# Row 4 ("Красное 1" office): shift code/name/type/parent one column left
$ws.Range("A4").Value = "'231101"
$ws.Range("B4").Value = "Красное 1"
$ws.Range("C4").Value = "ОПС"
$ws.Range("D4").Value = "'231100"

# The old "id" values (row 3/4 column A) and the now-redundant trailing
# parentOfficeCode column (E) are no longer part of the table
$ws.Range("E2:E4").Clear()

# Column B (office name) is widest now that it holds the name column
$ws.Columns.Item(2).ColumnWidth = 25.85

# Match the saved selection/cursor position
$ws.Range("F8").Select()
